$wb = $excel.ActiveWorkbook

# --- Add the second worksheet (工作表2) right after 工作表1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "工作表2"

# --- 工作表1: extend header row and existing rows with new columns ---
$ws1.Range("C1").Value = "std_dep"
$ws1.Range("D1").Value = "std_tel"
$ws1.Range("E1").Value = "櫃台"

$ws1.Range("C2").Value = "資工"
$ws1.Range("C3").Value = "國企"
$ws1.Range("C4").Value = "化材"

$ws1.Range("A4").Value = 0
$ws1.Range("B4").Value = "keke"

$ws1.Range("D2").Formula = '=VLOOKUP(C2,工作表2!$A$2:$B$5,2,0)'
$ws1.Range("D3").Formula = '=VLOOKUP(C3,工作表2!$A$2:$B$5,2,0)'
$ws1.Range("D4").Formula = '=VLOOKUP(C4,工作表2!$A$2:$B$5,2,0)'

$ws1.Range("E2").Formula = '=VLOOKUP(C2,工作表2!$A$1:$CL$5,3,0)'
$ws1.Range("E3").Formula = '=VLOOKUP(C3,工作表2!$A$1:$CL$5,3,0)'
$ws1.Range("E4").Formula = '=VLOOKUP(C4,工作表2!$A$1:$CL$5,3,0)'

# --- 工作表2: department / extension / counter lookup table ---
$ws2.Range("A1").Value = "科系"
$ws2.Range("B1").Value = "分機"
$ws2.Range("C1").Value = "櫃台"

$ws2.Range("A2").Value = "資工"
$ws2.Range("B2").Value = 12123
$ws2.Range("C2").Value = 1

$ws2.Range("A3").Value = "化材"
$ws2.Range("B3").Value = 121212
$ws2.Range("C3").Value = 2

$ws2.Range("A4").Value = "金資"
$ws2.Range("B4").Value = 12121212
$ws2.Range("C4").Value = 3

$ws2.Range("A5").Value = "國企"
$ws2.Range("B5").Value = 124887
$ws2.Range("C5").Value = 4

# --- Selections to mirror the authored state ---
[void]$ws2.Range("F7").Select()
[void]$ws1.Select()
[void]$ws1.Range("Q11").Select()
